# Daily attendance processing - 2026-01-04 19:07:16
# Swap the order of "System" and the recorder's email in column G
# (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
